# "cosmetic changes in the main script"
# The document title currently reads:
#   "The effect of parental presence on amygdala and mPFC activation
#    during fear conditioning: An exploratory study"
# It should read "effects" (plural) instead of "effect":
#   "The effects of parental presence on amygdala and mPFC activation
#    during fear conditioning: An exploratory study"
#
# Use Find/Replace (scoped to the title text) so the run keeps its
# existing run-level formatting (bold, white highlight, eastAsia font).

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute("The effect of parental presence on amygdala and ", `
              $true, `
              $false, `
              $false, `
              $false, `
              $false, `
              $true, `
              1, `
              $false, `
              "The effects of parental presence on amygdala and ", `
              2)
